# Update countries & provincias Spain
# Refresh the COVID-19 stats snapshot: update the "last updated" timestamp,
# update case/death/recovery numbers for the countries whose figures moved,
# and re-sort the 3 adjacent-row pairs whose ranking flipped as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 23:04"

# --- Rows whose ranking order swapped with their neighbour ---
# Surinam (was row 121) <-> Cabo Verde (was row 122): Cabo Verde overtakes Surinam
$ws.Range("A121").Value = "Cabo Verde"
$ws.Range("B121").Value = 4400
$ws.Range("C121").Value = 42
$ws.Range("D121").Value = 3851
$ws.Range("E121").Value = 507
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 42

$ws.Range("A122").Value = "Surinam"
$ws.Range("B122").Value = 4360
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 3544
$ws.Range("E122").Value = 725
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 91

# Republica de Chipre (was row 155) <-> Togo (was row 156): Togo overtakes Republica de Chipre
$ws.Range("A155").Value = "Togo"
$ws.Range("B155").Value = 1513
$ws.Range("C155").Value = 20
$ws.Range("D155").Value = 1127
$ws.Range("E155").Value = 352
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 34

$ws.Range("A156").Value = "Republica de Chipre"
$ws.Range("B156").Value = 1510
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 1237
$ws.Range("E156").Value = 251
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 22

# Isla de Man (was row 181) <-> Eritrea (was row 182): Eritrea overtakes Isla de Man
$ws.Range("A181").Value = "Eritrea"
$ws.Range("B181").Value = 341
$ws.Range("C181").Value = 11
$ws.Range("D181").Value = 295
$ws.Range("E181").Value = 46
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 0

$ws.Range("A182").Value = "Isla de Man"
$ws.Range("B182").Value = 337
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 312
$ws.Range("E182").Value = 1
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 24

# --- Remaining rows: numbers refreshed in place, same country/rank ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6503671
$ws.Range("C4").Value = 18096
$ws.Range("D4").Value = 3779728
$ws.Range("E4").Value = 2530087
$ws.Range("G4").Value = 322
$ws.Range("H4").Value = 193856

# Row 5: India
$ws.Range("B5").Value = 4367436
$ws.Range("C5").Value = 89852
$ws.Range("D5").Value = 3396027
$ws.Range("E5").Value = 897486
$ws.Range("G5").Value = 1107
$ws.Range("H5").Value = 73923

# Row 28: Israel
$ws.Range("B28").Value = 137565
$ws.Range("C28").Value = 3590
$ws.Range("D28").Value = 107003
$ws.Range("E28").Value = 29522

# Row 57: Costa Rica
$ws.Range("B57").Value = 49897
$ws.Range("C57").Value = 1117
$ws.Range("D57").Value = 19285
$ws.Range("E57").Value = 30081
$ws.Range("G57").Value = 21
$ws.Range("H57").Value = 531

# Row 68: Kenia
$ws.Range("B68").Value = 35356
$ws.Range("C68").Value = 151
$ws.Range("D68").Value = 21483
$ws.Range("E68").Value = 13274

# Row 82: Costa de Marfil
$ws.Range("B82").Value = 18778
$ws.Range("C82").Value = 77
$ws.Range("D82").Value = 17688
$ws.Range("E82").Value = 971

# Row 87: Senegal
$ws.Range("B87").Value = 14044
$ws.Range("C87").Value = 30
$ws.Range("E87").Value = 3715
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 292

# Row 116: Republica de Africa Central
$ws.Range("B116").Value = 4735
$ws.Range("C116").Value = 6
$ws.Range("D116").Value = 1825
$ws.Range("E116").Value = 2848

# Row 134: Angola
$ws.Range("B134").Value = 3033
$ws.Range("C134").Value = 52
$ws.Range("E134").Value = 1694
$ws.Range("G134").Value = 4
$ws.Range("H134").Value = 124

# Row 148: Sierra Leona
$ws.Range("B148").Value = 2064
$ws.Range("C148").Value = 9
$ws.Range("D148").Value = 1613
$ws.Range("E148").Value = 379
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = 72

# Row 180: Mauricio
$ws.Range("B180").Value = 361
$ws.Range("C180").Value = 5
$ws.Range("E180").Value = 16
